# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/value updates (safe from Excel numeric auto-conversion)
$ws.Cells.Item(2, 4).Value = "35.187.06"
$ws.Cells.Item(2, 5).Value = "  +1.08%  "
$ws.Cells.Item(3, 4).Value = "1.858.31"
$ws.Cells.Item(3, 5).Value = "  +1.58%  "
$ws.Cells.Item(4, 5).Value = "  +0.48%  "
$ws.Cells.Item(5, 5).Value = "  +3.75%  "
$ws.Cells.Item(6, 5).Value = "  +0.90%  "
$ws.Cells.Item(7, 5).Value = "  +0.41%  "
$ws.Cells.Item(8, 5).Value = "  +7.32%  "
$ws.Cells.Item(9, 5).Value = "  +0.88%  "
$ws.Cells.Item(10, 5).Value = "  +1.51%  "
$ws.Cells.Item(11, 5).Value = "  +0.17%  "
$ws.Cells.Item(12, 4).Value = "2.125.80"
$ws.Cells.Item(12, 5).Value = "  +1.53%  "
$ws.Cells.Item(13, 2).Value = "Chainlink"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(13, 5).Value = "  +1.63%  "
$ws.Cells.Item(14, 2).Value = "WrappedEther"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(14, 4).Value = "1.868.19"
$ws.Cells.Item(14, 5).Value = "  +2.28%  "
$ws.Cells.Item(15, 5).Value = "  +1.56%  "
$ws.Cells.Item(16, 5).Value = "  +1.72%  "
$ws.Cells.Item(17, 4).Value = "35.155.41"
$ws.Cells.Item(17, 5).Value = "  +0.96%  "
$ws.Cells.Item(19, 4).Value = "0.0₃0797"
$ws.Cells.Item(19, 5).Value = "  +1.32%  "
$ws.Cells.Item(20, 5).Value = "  +0.30%  "
$ws.Cells.Item(21, 5).Value = "  +0.56%  "
$ws.Cells.Item(22, 5).Value = "  +1.29%  "
$ws.Cells.Item(23, 5).Value = "  +0.44%  "
$ws.Cells.Item(24, 5).Value = "  -0.18%  "
$ws.Cells.Item(25, 5).Value = "  -1.97%  "
$ws.Cells.Item(26, 5).Value = "  +27.01%  "
$ws.Cells.Item(27, 5).Value = "  +3.29%  "
$ws.Cells.Item(28, 5).Value = "  +2.01%  "
$ws.Cells.Item(29, 5).Value = "  +0.42%  "
$ws.Cells.Item(30, 5).Value = "  +0.39%  "
$ws.Cells.Item(31, 5).Value = "  +1.30%  "
$ws.Cells.Item(33, 5).Value = "  +27.27%  "
$ws.Cells.Item(34, 5).Value = "  +2.41%  "
$ws.Cells.Item(35, 5).Value = "  +9.98%  "
$ws.Cells.Item(36, 5).Value = "  +17.21%  "
$ws.Cells.Item(37, 5).Value = "  +7.68%  "
$ws.Cells.Item(38, 5).Value = "  +5.57%  "
$ws.Cells.Item(39, 5).Value = "  +3.89%  "
$ws.Cells.Item(40, 5).Value = "  -1.57%  "
$ws.Cells.Item(41, 4).Value = "1.342.37"
$ws.Cells.Item(41, 5).Value = "  +0.05%  "
$ws.Cells.Item(42, 5).Value = "  +15.66%  "
$ws.Cells.Item(43, 5).Value = "  +3.99%  "
$ws.Cells.Item(44, 5).Value = "  +3.21%  "
$ws.Cells.Item(45, 5).Value = "  -0.09%  "
$ws.Cells.Item(46, 5).Value = "  +44.44%  "
$ws.Cells.Item(47, 5).Value = "  -0.50%  "
$ws.Cells.Item(48, 5).Value = "  +5.12%  "
$ws.Cells.Item(49, 4).Value = "2.038.69"
$ws.Cells.Item(49, 5).Value = "  +1.45%  "
$ws.Cells.Item(50, 5).Value = "  +1.08%  "
$ws.Cells.Item(51, 5).Value = "  +0.43%  "

# Numeric-looking text values: force text format so Excel keeps them as strings,
# then restore the default "Normal" style so no stray formatting is introduced.
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "239.42"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.623"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = "42.19"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0694"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0989"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "11.50"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.677"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "69.87"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "240.53"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "12.23"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "168.57"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.91"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.99"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "17.65"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.32"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0200"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "89.96"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0602"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "14.97"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.32"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "12.32"
$cell.Style = "Normal"
